$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.445.58'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.428.72'
$ws.Range("E3").Value = '  +6.79%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '295.89'
$ws.Range("E5").Value = '  -1.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.90'
$ws.Range("E6").Value = '  -2.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("E7").Value = '  +0.73%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.511'
$ws.Range("E9").Value = '  +0.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.47'
$ws.Range("E10").Value = '  +0.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0786'
$ws.Range("E11").Value = '  -1.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.14'
$ws.Range("E12").Value = '  +0.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  +1.96%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.793.99'
$ws.Range("E14").Value = '  +6.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.424.66'
$ws.Range("E15").Value = '  +6.73%  '

$ws.Range("E16").Value = '  +6.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.06'
$ws.Range("E17").Value = '  +2.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '46.287.10'
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.73'
$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0949'
$ws.Range("E20").Value = '  -2.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.25'
$ws.Range("E21").Value = '  +6.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.54'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '245.37'
$ws.Range("E23").Value = '  -1.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.80'
$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("E25").Value = '  +4.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '39.81'
$ws.Range("E27").Value = '  -3.92%  '

$ws.Range("E28").Value = '  +0.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.79'

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.42'
$ws.Range("E30").Value = '  +5.85%  '

$ws.Range("B31").Value = 'LidoDAOToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.81'
$ws.Range("E31").Value = '  +11.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.74'
$ws.Range("E32").Value = '  -2.41%  '

$ws.Range("E33").Value = '  +3.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '147.81'
$ws.Range("E34").Value = '  +0.15%  '

$ws.Range("E35").Value = '  -0.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.98'
$ws.Range("E36").Value = '  +16.79%  '

$ws.Range("E37").Value = '  +0.75%  '

$ws.Range("E38").Value = '  +0.17%  '

$ws.Range("E39").Value = '  -2.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.93'
$ws.Range("E40").Value = '  +1.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0302'
$ws.Range("E41").Value = '  +1.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.25'
$ws.Range("E42").Value = '  +3.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.972.18'
$ws.Range("E43").Value = '  +9.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.99'
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.83'
$ws.Range("E46").Value = '  -2.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.46'
$ws.Range("E47").Value = '  +31.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.60'
$ws.Range("E48").Value = '  +7.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.19'
$ws.Range("E49").Value = '  +6.49%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.666.05'
$ws.Range("E50").Value = '  +6.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.185'
$ws.Range("E51").Value = '  -0.54%  '

